# Apply the 0.1.6* -> 0.2.0 changelog update to the "Source table" sheet.
# This mirrors the commit: "list of functionalities find/replace 0.1.6* -> 0.2.0"
# plus the related note text tweaks (markdown link syntax fix, wording update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source table")
$ws.Activate()

# Replace every "version supported since" cell that still reads the old
# pre-release placeholder "0.1.6*" with the released version "0.2.0".
$ws.Range("D3").Value = "0.2.0"
$ws.Range("D22").Value = "0.2.0"

# Update the accompanying notes: bump the referenced version number and
# switch the bugfix reference to proper markdown link syntax.
$ws.Range("G31").Value = "Moved to io.crosssections in 0.2.0"
$ws.Range("G59").Value = "Moved to io.rr in 0.2.0"
$ws.Range("G60").Value = "Moved to io.rr in 0.2.0"
$ws.Range("G3").Value = "Critical bugfix for [#127](https://github.com/Deltares/HYDROLIB-core/issues/127)."

$ws.Range("D40").Value = "0.2.0"
$ws.Range("D41").Value = "0.2.0"
$ws.Range("D43").Value = "0.2.0"
$ws.Range("D44").Value = "0.2.0"
$ws.Range("D62").Value = "0.2.0"
$ws.Range("D63").Value = "0.2.0"

# Leave the selection on G4, matching the author's last interaction with
# the sheet (and reset the scroll position back to the top).
$ws.Range("A1").Select()
$ws.Range("G4").Select()
